$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.369.97'
$ws.Range('E2').Value = '  +1.98%  '
$ws.Range('D3').Value = '2.641.93'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '605.75'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '151.85'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.390'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +7.93%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.69'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.83'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').Value = '3.118.16'
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '64.218.03'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000150'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.43%  '
$ws.Range('D17').Value = '2.645.23'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.26'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +8.60%  '
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '353.36'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('E21').Value = '  +1.69%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('E23').Value = '  +3.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.85'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.75'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +13.58%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.73'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.67%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.40'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +8.66%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '548.99'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  +1.85%  '
$ws.Range('D33').Value = '0.0₃0873'
$ws.Range('E33').Value = '  +8.71%  '
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.33'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '167.75'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.92%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.04'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +8.95%  '
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +3.41%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '168.34'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '40.20'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.95'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.49%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0588'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.74%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '21.76'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('E47').Value = '  +16.05%  '
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('E49').Value = '  +3.14%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.43'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +5.31%  '
